$wb = $excel.ActiveWorkbook

# --- Layer0 sheet ---
$ws0 = $wb.Worksheets.Item("Layer0")

$ws0.Range("B2").Value = -0.8990630829450686
$ws0.Range("C2").Value = -0.6350191640325147

$ws0.Range("B3").Value = 0.4523604730746857
$ws0.Range("C3").Value = 0.5976472902761284

$ws0.Range("B4").Value = 1.231055233418578
$ws0.Range("C4").Value = -1.951909710443624

# --- Layer1 sheet ---
$ws1 = $wb.Worksheets.Item("Layer1")

$ws1.Range("B2").Value = -1.125485049979552
$ws1.Range("C2").Value = -0.2758537896692069

$ws1.Range("B3").Value = 1.011214613423473
$ws1.Range("C3").Value = 0.07082451861169821

$ws1.Range("B4").Value = -1.484071284767578
$ws1.Range("C4").Value = 0.9624889226984291
